$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B.
# This shifts the existing Unit column (B) -> C, Subcategories (C) -> D,
# and Tier (D) -> E, carrying along their values and styles.
[void]$ws.Columns("B:B").Insert()

# Match the width used for the new column in the target layout.
$ws.Columns("B:B").ColumnWidth = 12.5

# New "Value_type" column header and its first data rows.
$ws.Range("B1").Value = "Value_type"
$ws.Range("B2").Value = "Index"

# The (now shifted) "Unit" column used to hold "LCU"-prefixed labels;
# strip/replace them with the corresponding Value_type-relative labels.
$ws.Range("C2").Value = "% YoY, SA"
$ws.Range("C3").Value = "% YoY"
$ws.Range("C4").Value = "SA"
$ws.Range("C5").Value = "% MoM annualised"
$ws.Range("C6").Value = "% MoM annualised, SA"

# Update the active selection to match the onboarding instruction location.
[void]$ws.Range("C12").Select()
